$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: header labels for the new "how the algorithm works" block ---
# Order matches the original authoring sequence so shared-string ids line up.
$ws.Range("E14").Value = "Prize Percent "
$ws.Range("E14").Font.Bold = $true
$ws.Range("G14").Value = "Whole Numbers"
$ws.Range("G14").Font.Bold = $true
$ws.Range("J14").Value = "Round 1 Distribution"
$ws.Range("J14").Font.Bold = $true
$ws.Range("P14").Value = "Round 2 Distribution"
$ws.Range("P14").Font.Bold = $true
$ws.Range("S14").Value = "Acutal Percent"
$ws.Range("S14").Font.Bold = $true
$ws.Range("M14").Value = "Fractions"
$ws.Range("M14").Font.Bold = $true
$ws.Range("N14").Value = "Sorting"
$ws.Range("N14").Font.Bold = $true
$ws.Range("R14").Value = "TOTAL"
$ws.Range("R14").Font.Bold = $true

# --- Row 15 ---
$ws.Range("E15").Formula = "=C2"
$ws.Range("G15").Formula = "=FLOOR(C2,1)"
$ws.Range("J15").Formula = "=G15"
$ws.Range("M15").Formula = "=E15-G15"
$ws.Range("N15").Value = 2
$ws.Range("P15").Value = 1
$ws.Range("R15").Formula = "=J15+P15"
$ws.Range("S15").Formula = "=E2"

# --- Row 16 ---
$ws.Range("E16").Formula = "=C3"
$ws.Range("G16").Formula = "=FLOOR(C3,1)"
$ws.Range("J16").Formula = "=G16"
$ws.Range("M16").Formula = "=E16-G16"
$ws.Range("N16").Value = 4
$ws.Range("P16").Value = 0
$ws.Range("R16").Formula = "=J16+P16"
$ws.Range("S16").Formula = "=E3"

# --- Row 17 ---
$ws.Range("E17").Formula = "=C4"
$ws.Range("G17").Formula = "=FLOOR(C4,1)"
$ws.Range("J17").Formula = "=G17"
$ws.Range("M17").Formula = "=E17-G17"
$ws.Range("N17").Value = 6
$ws.Range("P17").Value = 0
$ws.Range("R17").Formula = "=J17+P17"
$ws.Range("S17").Formula = "=E4"

# --- Row 18 ---
$ws.Range("E18").Formula = "=C5"
$ws.Range("G18").Formula = "=FLOOR(C5,1)"
$ws.Range("J18").Formula = "=G18"
$ws.Range("M18").Formula = "=E18-G18"
$ws.Range("N18").Value = 5
$ws.Range("P18").Value = 0
$ws.Range("R18").Formula = "=J18+P18"
$ws.Range("S18").Formula = "=E5"

# --- Row 19 ---
$ws.Range("E19").Formula = "=C6"
$ws.Range("G19").Formula = "=FLOOR(C6,1)"
$ws.Range("J19").Formula = "=G19"
$ws.Range("M19").Formula = "=E19-G19"
$ws.Range("N19").Value = 1
$ws.Range("P19").Value = 1
$ws.Range("R19").Formula = "=J19+P19"
$ws.Range("S19").Formula = "=E6"

# --- Row 20 ---
$ws.Range("E20").Formula = "=C7"
$ws.Range("G20").Formula = "=FLOOR(C7,1)"
$ws.Range("J20").Formula = "=G20"
$ws.Range("M20").Formula = "=E20-G20"
$ws.Range("N20").Value = 3
$ws.Range("P20").Value = 1
$ws.Range("R20").Formula = "=J20+P20"
$ws.Range("S20").Formula = "=E7"

# --- Row 21: total check ---
$ws.Range("R21").Formula = "=SUM(R15:R20)"
$ws.Range("R21").Font.Bold = $true

# --- Column widths: Excel re-ran "best fit" on the columns touched by the
# new block (D/E split apart, G, J, M, N, P, S widened for the longer text) ---
$ws.Columns.Item(4).ColumnWidth = 13.608072916666666
$ws.Columns.Item(5).ColumnWidth = 11.498697916666666
$ws.Columns.Item(7).ColumnWidth = 13.830729166666666
$ws.Columns.Item(10).ColumnWidth = 17.608072916666668
$ws.Columns.Item(13).ColumnWidth = 7.721354166666667
$ws.Columns.Item(14).ColumnWidth = 6.166666666666667
$ws.Columns.Item(16).ColumnWidth = 17.608072916666668
$ws.Columns.Item(19).ColumnWidth = 12.498697916666666

# --- Selection moves to S14, matching the saved workbook state ---
$ws.Range("S14").Select()
